$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for each data row (2-50).
# Bump every value in C2:C50 from 46075 to 46076 (one day later),
# matching the source edit without touching any other cell/content.
for ($row = 2; $row -le 50; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value2 = 46076
    }
}
